# Generate Report for Handoff
# Updates the status/dates/error-detail for the file
# "41261206-05c9-4c23-85d7-253336fa2ef9.md" on the Overview, zh-cn and
# de-de sheets, and widens the "Error Detail" column to fit the new
# long message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff918ae04edf9c9fd30c3d5f3edd421539342cb5/e2e/41261206-05c9-4c23-85d7-253336fa2ef9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/208e0d847f330b9d9e2258924c19f1acb492e749/e2e/41261206-05c9-4c23-85d7-253336fa2ef9.md."

# --- Overview sheet ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 22:48:42"

# --- zh-cn sheet ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-24 22:48:37"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-24 22:48:42"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
